$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string (A1)
$ws.Range("A1").Value = "Datos actualizados a 29 de Marzo de 2020 a las 00:29"

# Re-point country names (column A) to reflect the reordered country list
$ws.Range("A68").Value = "Marruecos"
$ws.Range("A69").Value = "Lituania"
$ws.Range("A109").Value = "Nigeria"
$ws.Range("A110").Value = "Honduras"
$ws.Range("A111").Value = "Bielorrusia"
$ws.Range("A112").Value = "Martinica"
$ws.Range("A113").Value = "Camerun"
$ws.Range("A114").Value = "Georgia"
$ws.Range("A116").Value = "Trinidad yTobago"
$ws.Range("A117").Value = "Bolivia"

# Update statistic values (columns B-H) for affected rows
$ws.Range("B4").Value = 123271
$ws.Range("C4").Value = 19145
$ws.Range("E4").Value = 117838
$ws.Range("G4").Value = 506
$ws.Range("H4").Value = 2202
$ws.Range("E21").Value = 3784
$ws.Range("G21").Value = 22
$ws.Range("H21").Value = 114
$ws.Range("E35").Value = 1276
$ws.Range("G35").Value = 11
$ws.Range("H35").Value = 37
$ws.Range("B68").Value = 402
$ws.Range("C68").Value = 57
$ws.Range("D68").Value = 12
$ws.Range("E68").Value = 365
$ws.Range("F68").Value = 1
$ws.Range("H68").Value = 25
$ws.Range("B69").Value = 394
$ws.Range("C69").Value = 36
$ws.Range("D69").Value = 1
$ws.Range("E69").Value = 386
$ws.Range("F69").Value = 2
$ws.Range("H69").Value = 7
$ws.Range("B109").Value = 97
$ws.Range("E109").Value = 93
$ws.Range("F109").Value = 0
$ws.Range("B110").Value = 95
$ws.Range("C110").Value = 27
$ws.Range("D110").Value = 3
$ws.Range("E110").Value = 91
$ws.Range("F110").Value = 4
$ws.Range("H110").Value = 1
$ws.Range("B111").Value = 94
$ws.Range("D111").Value = 32
$ws.Range("E111").Value = 62
$ws.Range("F111").Value = 2
$ws.Range("H111").Value = 0
$ws.Range("B112").Value = 93
$ws.Range("D112").Value = 0
$ws.Range("E112").Value = 92
$ws.Range("F112").Value = 12
$ws.Range("H112").Value = 1
$ws.Range("B113").Value = 91
$ws.Range("C113").Value = 0
$ws.Range("D113").Value = 2
$ws.Range("E113").Value = 87
$ws.Range("F113").Value = 0
$ws.Range("H113").Value = 2
$ws.Range("B114").Value = 90
$ws.Range("C114").Value = 7
$ws.Range("D114").Value = 14
$ws.Range("E114").Value = 76
$ws.Range("F114").Value = 1
$ws.Range("H114").Value = 0
$ws.Range("B116").Value = 76
$ws.Range("C116").Value = 10
$ws.Range("D116").Value = 1
$ws.Range("E116").Value = 72
$ws.Range("G116").Value = 1
$ws.Range("H116").Value = 3
$ws.Range("C117").Value = 13
$ws.Range("D117").Value = 0
$ws.Range("E117").Value = 74
$ws.Range("G117").Value = 0
$ws.Range("H117").Value = 0
